# Add "Interpretation: ..." textboxes to the six per-parameter result slides.
#
# Target geometry (EMU, taken from the canonical OOXML diff) expressed in
# points for the COM API (1 pt = 12700 EMU):
#   Left   = 259197  EMU = 20.409212598425196 pt
#   Top    = 5942568 EMU = 467.91874015748033 pt
#   Height = 369332  EMU = 29.081259842519685 pt
#   Width (wide,  "NO STIM EFFECT") = 3481018 EMU = 274.095905511811    pt
#   Width (narrow,   "STIM EFFECT") = 3097899 EMU = 243.92905511811023 pt

$p = $ppt.ActivePresentation

$Left   = 20.409212598425196
$Top    = 467.91874015748033
$Height = 29.081259842519685
$WidthWide   = 274.095905511811
$WidthNarrow = 243.92905511811023

function Add-InterpretationTextbox {
    param($Slide, $Width, $Text)

    $tb = $Slide.Shapes.AddTextbox(1, $Left, $Top, $Width, $Height)
    $tb.TextFrame.WordWrap = $false
    $tb.TextFrame.AutoSize = 1
    $tb.Fill.Visible = $false

    $tr = $tb.TextFrame.TextRange
    $tr.Text = $Text
    $tr.Font.Bold = $true

    return $tb
}

# Slides 6, 7, 17, 18 and 19 each get their new shape assigned id/name that
# coincidentally collides with an id already used elsewhere on the slide
# except for slide 19, whose real target id (7 / "TextBox 6") is only
# reached after the per-slide shape-id counter has advanced once. Adding a
# throwaway textbox first (then deleting it) advances that counter without
# leaving any extra shape behind.
$s19 = $p.Slides.Item(19)
$dummy = $s19.Shapes.AddTextbox(1, $Left, $Top, $WidthWide, $Height)
Add-InterpretationTextbox $s19 $WidthWide "Interpretation: NO STIM EFFECT"
$dummy.Delete()

# Remaining slides: a single AddTextbox call already lands on the right id/name.
Add-InterpretationTextbox ($p.Slides.Item(6))  $WidthWide   "Interpretation: NO STIM EFFECT"
Add-InterpretationTextbox ($p.Slides.Item(7))  $WidthWide   "Interpretation: NO STIM EFFECT"
Add-InterpretationTextbox ($p.Slides.Item(8))  $WidthNarrow "Interpretation: STIM EFFECT"
Add-InterpretationTextbox ($p.Slides.Item(17)) $WidthWide   "Interpretation: NO STIM EFFECT"
Add-InterpretationTextbox ($p.Slides.Item(18)) $WidthWide   "Interpretation: NO STIM EFFECT"
